$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that currently hold the text "NULL" (shared string) and must become numeric 0.
$cells = @(
    "C83","D83",
    "C85","D85",
    "C104","D104",
    "C142","D142",
    "C150","D150",
    "B152","C152","D152",
    "C153","D153",
    "C164","D164",
    "C167","D167",
    "C186","D186",
    "C207","D207",
    "C212","D212",
    "C214","D214",
    "C220","D220",
    "C228","D228"
)

foreach ($c in $cells) {
    $ws.Range($c).Value = 0
}

# After replacing all "NULL" text occurrences, select the whole worksheet
# (mirrors the post Find&Replace-all selection state captured in the file).
$ws.Range("A1:XFD1048576").Select()
